$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("INVENTORY_TYPE", "STOCK_TYPE", "INVENTORY_STATUS", "SUBINVENTORIES", "COSTING_METHOD")

$row = 18
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row = $row + 1
}

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("A23").Select() | Out-Null
